$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data block (CE / HUGO LINO GONZALEZ MARTINEZ) inserted ahead of the
# existing CC / CRISTIAN RAMON FLOREZ ECHENIQUE rows, data refreshed with
# updated "Valor Mora" / "Salario Basico" amounts, per commit:
# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"

$rows = @(
    @{ Row = 16; B = "CE"; C = "386447";  D = "HUGO LINO GONZALEZ MARTINEZ";     E = "1809"; F = 40000; G = 1000000 },
    @{ Row = 17; B = "CE"; C = "386447";  D = "HUGO LINO GONZALEZ MARTINEZ";     E = "1810"; F = 40000; G = 1000000 },
    @{ Row = 18; B = "CE"; C = "386447";  D = "HUGO LINO GONZALEZ MARTINEZ";     E = "1812"; F = 40000; G = 1000000 },
    @{ Row = 19; B = "CC"; C = "9102460"; D = "CRISTIAN RAMON FLOREZ ECHENIQUE"; E = "1904"; F = 80000; G = 2000000 },
    @{ Row = 20; B = "CC"; C = "9102460"; D = "CRISTIAN RAMON FLOREZ ECHENIQUE"; E = "1906"; F = 80000; G = 2000000 },
    @{ Row = 21; B = "CC"; C = "9102460"; D = "CRISTIAN RAMON FLOREZ ECHENIQUE"; E = "1907"; F = 80000; G = 2000000 },
    @{ Row = 22; B = "CC"; C = "9102460"; D = "CRISTIAN RAMON FLOREZ ECHENIQUE"; E = "1908"; F = 80000; G = 2000000 },
    @{ Row = 23; B = "CC"; C = "9102460"; D = "CRISTIAN RAMON FLOREZ ECHENIQUE"; E = "1909"; F = 80000; G = 2000000 },
    @{ Row = 24; B = "CC"; C = "9102460"; D = "CRISTIAN RAMON FLOREZ ECHENIQUE"; E = "1910"; F = 80000; G = 2000000 },
    @{ Row = 25; B = "CC"; C = "9102460"; D = "CRISTIAN RAMON FLOREZ ECHENIQUE"; E = "1911"; F = 64000; G = 2000000 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
}
